# "Generate Report for Handback"
#
# The 5406c99d-7b31-4485-ba4c-18fecc54a078.md file has been handed back
# (its handback report was generated), so its status flips from
# "Ready for handoff" to "Handed back: in sync with en-US" and its rows
# move to the top of each sheet (sorted ahead of the other two, already
# handed-back files) on every sheet (Overview, zh-cn, de-de). The
# locale sheets also gain their "Latest Target File" / "Latest Handback
# File" / "Latest Handback DateTime" values now that the handback is in.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Clear out the existing hyperlinks on column B; they get rebuilt below
# once the rows are reordered.
$ws.Range("B2:B4").Hyperlinks.Delete()

# Row 2 becomes the 5406c99d entry (now handed back).
$ws.Range("A2").Value = "5406c99d-7b31-4485-ba4c-18fecc54a078.md"
$ws.Range("B2").Value = "e2e\5406c99d-7b31-4485-ba4c-18fecc54a078.md"
$ws.Range("C2").Value = ".md"
$ws.Range("E2").Value = "Handed back: in sync with en-US"
$ws.Range("F2").Value = "Handed back: in sync with en-US"
$ws.Range("G2").Value = "2016-08-28 13:04:19"

# Row 3 becomes the ffff5840b6e6 entry.
$ws.Range("A3").Value = "ffff5840b6e6-9309-4023-912d-e5b977376339.md"
$ws.Range("B3").Value = "e2e\ffff5840b6e6-9309-4023-912d-e5b977376339.md"
$ws.Range("C3").Value = ".md"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-08-28 13:03:19"

# Row 4 becomes the ffffff53b45ca6 entry.
$ws.Range("A4").Value = "ffffff53b45ca6-4475-403c-867e-cd8ab3e556ac.md"
$ws.Range("B4").Value = "e2e\ffffff53b45ca6-4475-403c-867e-cd8ab3e556ac.md"
$ws.Range("C4").Value = ".md"
$ws.Range("E4").Value = "Handed back: in sync with en-US"
$ws.Range("F4").Value = "Handed back: in sync with en-US"
$ws.Range("G4").Value = "2016-08-28 13:03:19"

$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/82d148b724c4c47760a7262eb3085e70a8254e7a/e2e/5406c99d-7b31-4485-ba4c-18fecc54a078.md", "", "", "e2e\5406c99d-7b31-4485-ba4c-18fecc54a078.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f4786dc9ee3a37b53ff10c392a1a3c637348e1d/e2e/ffff5840b6e6-9309-4023-912d-e5b977376339.md", "", "", "e2e\ffff5840b6e6-9309-4023-912d-e5b977376339.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/82d148b724c4c47760a7262eb3085e70a8254e7a/e2e/ffffff53b45ca6-4475-403c-867e-cd8ab3e556ac.md", "", "", "e2e\ffffff53b45ca6-4475-403c-867e-cd8ab3e556ac.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2:A4").Hyperlinks.Delete()
$ws.Range("I2:I4").Hyperlinks.Delete()

# Row 2 becomes the 5406c99d entry (now handed back, target + handback
# info now populated).
$ws.Range("A2").Value = "5406c99d-7b31-4485-ba4c-18fecc54a078.md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("F2").Value = "False"
$ws.Range("G2").Value = "5406c99d-7b31-4485-ba4c-18fecc54a078.dccf443320aed0861eec315d7149de779c4cefba.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-28 13:04:13"
$ws.Range("I2").Value = "5406c99d-7b31-4485-ba4c-18fecc54a078.md"
$ws.Range("J2").Value = "5406c99d-7b31-4485-ba4c-18fecc54a078.dccf443320aed0861eec315d7149de779c4cefba.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-28 13:04:43"

# Row 3 becomes the ffff5840b6e6 entry.
$ws.Range("A3").Value = "ffff5840b6e6-9309-4023-912d-e5b977376339.md"
$ws.Range("F3").Value = "False"
$ws.Range("G3").Value = "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.eeca0f1d8dd5ea7239f7a9106f95f68b32accc9f.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-28 13:03:14"
$ws.Range("I3").Value = "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.md"
$ws.Range("J3").Value = "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.eeca0f1d8dd5ea7239f7a9106f95f68b32accc9f.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-28 13:03:31"

# Row 4 becomes the ffffff53b45ca6 entry.
$ws.Range("A4").Value = "ffffff53b45ca6-4475-403c-867e-cd8ab3e556ac.md"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("F4").Value = "True"
$ws.Range("G4").Value = "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.eeca0f1d8dd5ea7239f7a9106f95f68b32accc9f.zh-cn.xlf"
$ws.Range("H4").Value = "2016-08-28 13:03:14"
$ws.Range("I4").Value = "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.md"
$ws.Range("J4").Value = "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.eeca0f1d8dd5ea7239f7a9106f95f68b32accc9f.zh-cn.xlf"
$ws.Range("K4").Value = "2016-08-28 13:03:31"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/82d148b724c4c47760a7262eb3085e70a8254e7a/e2e/5406c99d-7b31-4485-ba4c-18fecc54a078.md", "", "", "5406c99d-7b31-4485-ba4c-18fecc54a078.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/82d148b724c4c47760a7262eb3085e70a8254e7a/e2e/5406c99d-7b31-4485-ba4c-18fecc54a078.md", "", "", "5406c99d-7b31-4485-ba4c-18fecc54a078.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f4786dc9ee3a37b53ff10c392a1a3c637348e1d/e2e/ffff5840b6e6-9309-4023-912d-e5b977376339.md", "", "", "ffff5840b6e6-9309-4023-912d-e5b977376339.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5179297f75b2ccab043380b63aea4c3f4228d466/e2e/a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.md", "", "", "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/82d148b724c4c47760a7262eb3085e70a8254e7a/e2e/ffffff53b45ca6-4475-403c-867e-cd8ab3e556ac.md", "", "", "ffffff53b45ca6-4475-403c-867e-cd8ab3e556ac.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5179297f75b2ccab043380b63aea4c3f4228d466/e2e/a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.md", "", "", "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2:A4").Hyperlinks.Delete()
$ws.Range("I2:I4").Hyperlinks.Delete()

# Row 2 becomes the 5406c99d entry (now handed back, target + handback
# info now populated).
$ws.Range("A2").Value = "5406c99d-7b31-4485-ba4c-18fecc54a078.md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("F2").Value = "False"
$ws.Range("G2").Value = "5406c99d-7b31-4485-ba4c-18fecc54a078.dccf443320aed0861eec315d7149de779c4cefba.de-de.xlf"
$ws.Range("H2").Value = "2016-08-28 13:04:19"
$ws.Range("I2").Value = "5406c99d-7b31-4485-ba4c-18fecc54a078.md"
$ws.Range("J2").Value = "5406c99d-7b31-4485-ba4c-18fecc54a078.dccf443320aed0861eec315d7149de779c4cefba.de-de.xlf"
$ws.Range("K2").Value = "2016-08-28 13:04:50"

# Row 3 becomes the ffff5840b6e6 entry.
$ws.Range("A3").Value = "ffff5840b6e6-9309-4023-912d-e5b977376339.md"
$ws.Range("F3").Value = "False"
$ws.Range("G3").Value = "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.eeca0f1d8dd5ea7239f7a9106f95f68b32accc9f.de-de.xlf"
$ws.Range("H3").Value = "2016-08-28 13:03:19"
$ws.Range("I3").Value = "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.md"
$ws.Range("J3").Value = "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.eeca0f1d8dd5ea7239f7a9106f95f68b32accc9f.de-de.xlf"
$ws.Range("K3").Value = "2016-08-28 13:03:38"

# Row 4 becomes the ffffff53b45ca6 entry.
$ws.Range("A4").Value = "ffffff53b45ca6-4475-403c-867e-cd8ab3e556ac.md"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("F4").Value = "True"
$ws.Range("G4").Value = "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.eeca0f1d8dd5ea7239f7a9106f95f68b32accc9f.de-de.xlf"
$ws.Range("H4").Value = "2016-08-28 13:03:19"
$ws.Range("I4").Value = "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.md"
$ws.Range("J4").Value = "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.eeca0f1d8dd5ea7239f7a9106f95f68b32accc9f.de-de.xlf"
$ws.Range("K4").Value = "2016-08-28 13:03:38"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/82d148b724c4c47760a7262eb3085e70a8254e7a/e2e/5406c99d-7b31-4485-ba4c-18fecc54a078.md", "", "", "5406c99d-7b31-4485-ba4c-18fecc54a078.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/82d148b724c4c47760a7262eb3085e70a8254e7a/e2e/5406c99d-7b31-4485-ba4c-18fecc54a078.md", "", "", "5406c99d-7b31-4485-ba4c-18fecc54a078.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f4786dc9ee3a37b53ff10c392a1a3c637348e1d/e2e/ffff5840b6e6-9309-4023-912d-e5b977376339.md", "", "", "ffff5840b6e6-9309-4023-912d-e5b977376339.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d84cb041d1a9ed21083c1d0d2874a41ced3196cd/e2e/a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.md", "", "", "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/82d148b724c4c47760a7262eb3085e70a8254e7a/e2e/ffffff53b45ca6-4475-403c-867e-cd8ab3e556ac.md", "", "", "ffffff53b45ca6-4475-403c-867e-cd8ab3e556ac.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d84cb041d1a9ed21083c1d0d2874a41ced3196cd/e2e/a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.md", "", "", "a3df63f6-f3ab-4dbf-818e-7fbc60ab89c2.md") | Out-Null

$ws = $wb.Worksheets.Item("Overview")
$ws.Select()
$ws.Range("A1").Select()
